{"js": "// The author typed a leading space at the very start of the document, which is\n// why Word's \"last edit\" bookmark (_GoBack) moved from the end of the document\n// (after \"Vraag ver\" + \"d\") to the beginning (right after the new leading space).\n// We reproduce both effects: insert the leading space, move the _GoBack bookmark,\n// and normalize the run that used to be split around the old bookmark location.\n\nconst body = context.document.body;\n\n// 1) Remove the existing \"_GoBack\" bookmark (currently sits inside the last\n//    paragraph, between \"Vraag ver\" + \"d\" and \"er aan Casper...\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Re-join the run that the old bookmark used to split (\"d\" | \"er aan Casper...\")\n//    back into a single run by re-typing its combined text over itself.\nconst tailResults = body.search(\n  \"der aan Casper als hij opmerkingen heeft over wat misschien makkelijk beter kan \",\n  { matchCase: true }\n);\ntailResults.load(\"items\");\nawait context.sync();\n\nif (tailResults.items.length > 0) {\n  tailResults.items[0].insertText(\n    \"der aan Casper als hij opmerkingen heeft over wat misschien makkelijk beter kan \",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 3) Insert a leading space at the very start of the document body.\nconst startRange = body.getRange(\"Start\");\nstartRange.insertText(\" \", \"Before\");\nawait context.sync();\n\n// 4) Re-create the \"_GoBack\" bookmark collapsed right after that new leading\n//    space (i.e. immediately before \"Wat moet er nog gebeuren...\").\nconst headResults = body.search(\"Wat moet er nog gebeuren\", { matchCase: true });\nheadResults.load(\"items\");\nawait context.sync();\n\nif (headResults.items.length > 0) {\n  const pointBeforeWat = headResults.items[0].getRange(\"Start\");\n  pointBeforeWat.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author typed a leading space at the very start of the document, which is\n# why Word's \"last edit\" bookmark (_GoBack) moved from the end of the document\n# (between \"Vraag verd\" and \"er aan Casper...\") to the beginning (right after\n# the new leading space). We reproduce both effects below.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# Step 1: protect the \"Vraag ve\" / \"r\" run boundary with a temporary bookmark\n#         so later edits cannot merge those two runs together. (In this Word\n#         engine a bookmark acts as a hard barrier against run-merging.)\n# ---------------------------------------------------------------------------\n$rBoundary = $d.Content\n$rBoundary.Find.ClearFormatting()\n$foundBoundary = $rBoundary.Find.Execute(\"Vraag ver\")\nif ($foundBoundary) {\n    $barrierRange = $d.Range($rBoundary.End, $rBoundary.End)\n    $d.Bookmarks.Add(\"TEMP_BARRIER\", $barrierRange)\n}\n\n# ---------------------------------------------------------------------------\n# Step 2: remove the existing \"_GoBack\" bookmark. It currently sits right\n#         between the old \"d\" run and the \"er aan Casper...\" run.\n# ---------------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# ---------------------------------------------------------------------------\n# Step 3: touch the text inside \"er aan Casper...\" (insert then remove a\n#         placeholder) so Word re-joins that run with the now-adjacent \"d\"\n#         run into a single run, bounded on the left by TEMP_BARRIER so\n#         \"Vraag ve\"/\"r\" stay untouched.\n# ---------------------------------------------------------------------------\n$rCasper = $d.Content\n$rCasper.Find.ClearFormatting()\n$foundCasper = $rCasper.Find.Execute(\"Casper\")\nif ($foundCasper) {\n    $rTouch = $d.Range($rCasper.Start, $rCasper.Start)\n    $rTouch.InsertBefore(\"TEMP\")\n}\n\n$rCleanup = $d.Content\n$rCleanup.Find.ClearFormatting()\n$foundCleanup = $rCleanup.Find.Execute(\"TEMP\")\nif ($foundCleanup) {\n    $rCleanup.Text = \"\"\n}\n\n# ---------------------------------------------------------------------------\n# Step 4: remove the temporary barrier bookmark.\n# ---------------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"TEMP_BARRIER\")) {\n    $d.Bookmarks.Item(\"TEMP_BARRIER\").Delete()\n}\n\n# ---------------------------------------------------------------------------\n# Step 5: insert a leading space at the very start of the document body.\n# ---------------------------------------------------------------------------\n$rStart = $d.Range(0, 0)\n$rStart.InsertBefore(\" \")\n\n# ---------------------------------------------------------------------------\n# Step 6: re-create the \"_GoBack\" bookmark collapsed right after that new\n#         leading space (i.e. immediately before \"Wat moet er nog gebeuren...\").\n# ---------------------------------------------------------------------------\n$rHead = $d.Content\n$rHead.Find.ClearFormatting()\n$foundHead = $rHead.Find.Execute(\"Wat moet er nog gebeuren\")\nif ($foundHead) {\n    $bmRange = $d.Range($rHead.Start, $rHead.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
